$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 0.0150041999
$ws.Range("F5").Value = 0.0004644361
$ws.Range("G5").Value = 0.1127824783
$ws.Range("H5").Value = 5.2422862053
$ws.Range("I5").Value = 0.0193367973
$ws.Range("J5").Value = 0.0009260924
$ws.Range("K5").Value = 0.1941526532
$ws.Range("L5").Value = 14.5507297516
$ws.Range("M5").Value = 0.08713245520243898
$ws.Range("N5").Value = 0.01322879565028934
$ws.Range("O5").Value = 0.3662151112999998
$ws.Range("P5").Value = 0.01917805240263791
$ws.Range("Q5").Value = 0.2049857087914634
$ws.Range("R5").Value = 0.04998100075705328
$ws.Range("S5").Value = 0.4581298399999998
$ws.Range("T5").Value = 0.04484785861462664
$ws.Range("U5").Value = 0.09358293865731712
$ws.Range("V5").Value = 0.01512463012686926
$ws.Range("W5").Value = 0.3930178176000005
$ws.Range("X5").Value = 0.01288598188120367
$ws.Range("Y5").Value = 0.1571766717621952
$ws.Range("Z5").Value = 0.03214146853309562
$ws.Range("AA5").Value = 0.4359028488999996
$ws.Range("AB5").Value = 0.02223979184531321
$ws.Range("AC5").Value = 0.08520624953354147
$ws.Range("AD5").Value = 0.01082621235206209
$ws.Range("AE5").Value = 0.4116257720999998
$ws.Range("AF5").Value = 0.01730429481936817
$ws.Range("AG5").Value = 0.215171872347636
$ws.Range("AH5").Value = 0.05024528974572511
$ws.Range("AI5").Value = 0.5421935017999999
$ws.Range("AJ5").Value = 0.04373778559697056

$ws.Range("E6").Value = 0.0134369098
$ws.Range("F6").Value = 0.0003546652
$ws.Range("G6").Value = 0.078820318
$ws.Range("H6").Value = 4.3206686974
$ws.Range("I6").Value = 0.0184779037
$ws.Range("J6").Value = 0.0007264043
$ws.Range("K6").Value = 0.1395274401
$ws.Range("L6").Value = 15.690653801
$ws.Range("M6").Value = 0.09365449193048778
$ws.Range("N6").Value = 0.01364614231864257
$ws.Range("O6").Value = 0.3431521409
$ws.Range("P6").Value = 0.02095367475148259
$ws.Range("Q6").Value = 0.2871720668256098
$ws.Range("R6").Value = 0.09924895117487342
$ws.Range("S6").Value = 0.5757091534000001
$ws.Range("T6").Value = 0.06378066916460212
$ws.Range("U6").Value = 0.08409094357682924
$ws.Range("V6").Value = 0.01235060334705156
$ws.Range("W6").Value = 0.3420374807000002
$ws.Range("X6").Value = 0.01180574616415395
$ws.Range("Y6").Value = 0.1384540965536586
$ws.Range("Z6").Value = 0.03034102997371167
$ws.Range("AA6").Value = 0.4383362913000006
$ws.Range("AB6").Value = 0.02168252878731586
$ws.Range("AC6").Value = 0.08318048070847457
$ws.Range("AD6").Value = 0.01047380758161532
$ws.Range("AE6").Value = 0.4135514535000002
$ws.Range("AF6").Value = 0.01693586533446379
$ws.Range("AG6").Value = 0.4277105660935772
$ws.Range("AH6").Value = 0.1906173135636298
$ws.Range("AI6").Value = 0.7572189151000002
$ws.Range("AJ6").Value = 0.08721984374706636
